$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the previously-empty row 50 with new sample entries
$ws.Range("A50").Value = "Inj"
$ws.Range("B50").Value = "Stock"

# Reflect the final active selection from the edit session
$ws.Range("B7").Select()
